$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new work-log entry for the latest meeting
$ws.Cells.Item(35, 2).Value = "Meeting - Implementing Multiline and MultipleChoice questions in Form edit view"
$ws.Cells.Item(35, 3).Value = 1
$ws.Cells.Item(35, 4).Value = 40886

# Match the date formatting used by the row above (column D is a date column)
$ws.Range("D34").Copy()
$ws.Range("D35").PasteSpecial(-4122)

# Move the active selection down, as Excel does after finishing data entry
$ws.Range("B36").Select() | Out-Null
